# "Update metodo de la ingenieria"
#
# Content changes:
#  - Sheet "Metodo" (sheet3): clear D3 ("Revisar en el postmortem" note), which
#    removes that now-unused shared string, and check off ("x") the remaining
#    "Fase 3/4/5" rows (C12, C14, C15, C17, C18, C20).
#  - Sheet "Entregables" (sheet1): check off ("x") C14, C15, C16 (Fase 3/4/5
#    under "Informe metodo de la ingenieria") and C22 ("Actualizar el diagrama").
#  - Sheet "Rubrica" (sheet2): check off ("x") C3, C11, C12.
#  - Sheet "Conexion a interfaz" (sheet4): no content changes.
#
# Selection/view-state is updated to match too (best effort - the last
# worksheet whose Range.Select() is called ends up as the active tab, so
# sheet4 is selected last to keep it the active sheet, matching the source
# workbook where activeTab stays pointed at "Conexion a interfaz").

$wb = $excel.ActiveWorkbook

$wsEntregables = $wb.Worksheets.Item(1)
$wsRubrica     = $wb.Worksheets.Item(2)
$wsMetodo      = $wb.Worksheets.Item(3)
$wsConexion    = $wb.Worksheets.Item(4)

# --- Sheet "Entregables" ---------------------------------------------------
$wsEntregables.Range("C14").Value = "x"
$wsEntregables.Range("C15").Value = "x"
$wsEntregables.Range("C16").Value = "x"
$wsEntregables.Range("C22").Value = "x"

# --- Sheet "Rubrica" ---------------------------------------------------------
$wsRubrica.Range("C3").Value  = "x"
$wsRubrica.Range("C11").Value = "x"
$wsRubrica.Range("C12").Value = "x"

# --- Sheet "Metodo" ----------------------------------------------------------
$wsMetodo.Range("D3").ClearContents()
$wsMetodo.Range("C12").Value = "x"
$wsMetodo.Range("C14").Value = "x"
$wsMetodo.Range("C15").Value = "x"
$wsMetodo.Range("C17").Value = "x"
$wsMetodo.Range("C18").Value = "x"
$wsMetodo.Range("C20").Value = "x"

# --- Sheet "Conexion a interfaz" --------------------------------------------
# (no cell content changes)

# --- Selections (best effort; keeps "Conexion a interfaz" as the active tab)
$wsEntregables.Range("C17").Select()
$wsRubrica.Range("C13").Select()
$wsMetodo.Range("E10").Select()
$wsConexion.Range("D5").Select()
